# Trade #16 closed at 2026-02-17 07:59:13 - unknown UNKNOWN +0.000%
$wb = $excel.ActiveWorkbook

# ---- Summary sheet ----
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B5").Value = -0.1      # Total P&L %
$wsSummary.Range("B6").Value = 16        # Total Trades
$wsSummary.Range("B9").Value = 31.25     # Win Rate %

# ---- Strategy Status sheet ----
$wsStatus = $wb.Worksheets.Item("Strategy Status")
$wsStatus.Range("D4").Value = 16         # Trades (MarketMaking)
$wsStatus.Range("G4").Value = 31.25      # Win Rate % (MarketMaking)

# ---- Helper to append the new trade row (row 17) to a trades sheet ----
function Add-Trade16Row($ws) {
    $row = 17

    $ws.Cells.Item($row, 1).Value = 16

    $ws.Cells.Item($row, 2).NumberFormat = "@"
    $ws.Cells.Item($row, 2).Value = "2026-02-17"

    $ws.Cells.Item($row, 3).NumberFormat = "@"
    $ws.Cells.Item($row, 3).Value = "07:59:07"

    $ws.Cells.Item($row, 4).NumberFormat = "@"
    $ws.Cells.Item($row, 4).Value = "MarketMaking"

    $ws.Cells.Item($row, 5).NumberFormat = "@"
    $ws.Cells.Item($row, 5).Value = "UP"

    $ws.Cells.Item($row, 6).Value = 0.03
    $ws.Cells.Item($row, 7).Value = 0.03

    $ws.Cells.Item($row, 8).NumberFormat = "@"
    $ws.Cells.Item($row, 8).Value = "CLOSED"

    $ws.Cells.Item($row, 9).Value = 0
    $ws.Cells.Item($row, 10).Value = 0
    $ws.Cells.Item($row, 11).Value = 99.92
    $ws.Cells.Item($row, 12).Value = 0
    $ws.Cells.Item($row, 13).Value = 0
    $ws.Cells.Item($row, 14).Value = 0.6

    $ws.Cells.Item($row, 15).NumberFormat = "@"
    $ws.Cells.Item($row, 15).Value = "Normal spread capture: 19600 bps"

    $ws.Cells.Item($row, 16).NumberFormat = "@"
    $ws.Cells.Item($row, 16).Value = "early_exit"

    $ws.Cells.Item($row, 17).Value = 0.14
}

# ---- All Trades sheet ----
$wsAllTrades = $wb.Worksheets.Item("All Trades")
Add-Trade16Row $wsAllTrades

# ---- MarketMaking sheet ----
$wsMarketMaking = $wb.Worksheets.Item("MarketMaking")
Add-Trade16Row $wsMarketMaking
